# Add a new "3d_classic" worksheet right after the existing "3d" sheet.
# It holds the same data as "3d" but laid out in "classic" (stacked / long)
# format: one row per (a, b) combination, with columns a, b, c0, c1, c2.

$wb = $excel.ActiveWorkbook
$ws3d = $wb.Worksheets.Item("3d")

$newWs = $wb.Worksheets.Add($null, $ws3d)
$newWs.Name = "3d_classic"

$data = @(
    @("a", "b",  "c0", "c1", "c2"),
    @(1,   "b0", 0,    1,    2),
    @(1,   "b1", 3,    4,    5),
    @(2,   "b0", 6,    7,    8),
    @(2,   "b1", 9,    10,   11),
    @(3,   "b0", 12,   13,   14),
    @(3,   "b1", 15,   16,   17)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $newWs.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
